$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Widen column A (19.85546875 -> 26.7109375 characters)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 26.7109375

# ---------------------------------------------------------------------------
# 2. Add hyperlinks first (so the banding style applied afterwards is not
#    overwritten by the default hyperlink style Excel applies automatically)
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("A50"), "https://www.moddb.com/mods/stalker-anomaly/addons/headlamp-animation-fix") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A51"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A52"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A53"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A55"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A54"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A57"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A59"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A56"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A58"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A60"), "https://www.moddb.com/mods/stalker-anomaly/addons/utjans-item-ui-improvements") | Out-Null

# ---------------------------------------------------------------------------
# 3. Re-apply the correct alternating row banding/style by copying it over
#    from the existing rows that already carry the right look (this also
#    restores the style that the Hyperlinks.Add call above just overwrote).
#    Row 48 = "even" banding (s=12,19,11,9,15), Row 49 = "odd" banding
#    (s=13,20,10,3,14).
# ---------------------------------------------------------------------------
$ws.Range("A48:E48").Copy($ws.Range("A50:E50"))   # even banding
$ws.Range("A49:E49").Copy($ws.Range("A51:E51"))   # odd banding
$ws.Range("A48:E48").Copy($ws.Range("A52:E52"))   # even banding
$ws.Range("A49:E49").Copy($ws.Range("A53:E53"))   # odd banding
$ws.Range("A48:E48").Copy($ws.Range("A54:E54"))   # even banding
$ws.Range("A49:E49").Copy($ws.Range("A55:E55"))   # odd banding
$ws.Range("A48:E48").Copy($ws.Range("A56:E56"))   # even banding
$ws.Range("A49:E49").Copy($ws.Range("A57:E57"))   # odd banding
$ws.Range("A48:E48").Copy($ws.Range("A58:E58"))   # even banding
$ws.Range("A49:E49").Copy($ws.Range("A59:E59"))   # odd banding
$ws.Range("A48:E48").Copy($ws.Range("A60:E60"))   # even banding

# Row 61 stays empty, but its "A" cell style switches from the old unfilled
# style to the hyperlink-banded "odd" style (matching its siblings), ready
# for a future entry.
$ws.Range("A49").Copy($ws.Range("A61"))
$ws.Range("A61").ClearContents()

# ---------------------------------------------------------------------------
# 4. Write the actual cell values for the two newly tracked mods.
# ---------------------------------------------------------------------------

# Row 50 - Headlamp anim fix (Ravenascendant)
$ws.Range("A50").Value2 = "Headlamp anim fix"
$ws.Range("B50").Value2 = "zz_headlamp_animation_fix.script"
$ws.Range("C50").Value2 = "scripts/"
$ws.Range("D50").Value2 = "Yes"
$ws.Range("E50").Value2 = "Ravenascendant"

# Rows 51-60 - Utjans Item Ui Improvements (Utjan)
$ws.Range("A51").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B51").Value2 = "eq_icons_settings.ltx"
$ws.Range("C51").Value2 = "configs/plugins/"
$ws.Range("D51").Value2 = "Yes"
$ws.Range("E51").Value2 = "Utjan"

$ws.Range("A52").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B52").Value2 = "ui_st_item_icon_info.xml"
$ws.Range("C52").Value2 = "configs/text/eng/"
$ws.Range("D52").Value2 = "Yes"
$ws.Range("E52").Value2 = "Utjan"

$ws.Range("A53").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B53").Value2 = "ui_st_item_icon_info.xml"
$ws.Range("C53").Value2 = "configs/text/rus/"
$ws.Range("D53").Value2 = "Yes"
$ws.Range("E53").Value2 = "Utjan"

$ws.Range("A54").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B54").Value2 = "ui_dyn_eq_icons.xml"
$ws.Range("C54").Value2 = "configs/ui/textures_descr"
$ws.Range("D54").Value2 = "Yes"
$ws.Range("E54").Value2 = "Utjan"

$ws.Range("A55").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B55").Value2 = "ui_utjan_icon_info.xml"
$ws.Range("C55").Value2 = "configs/ui/"
$ws.Range("D55").Value2 = "Yes"
$ws.Range("E55").Value2 = "Utjan"

$ws.Range("A56").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B56").Value2 = "dynamic_eq_icon.script"
$ws.Range("C56").Value2 = "scripts/"
$ws.Range("D56").Value2 = "Yes"
$ws.Range("E56").Value2 = "Utjan"

$ws.Range("A57").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B57").Value2 = "item_icon_info_mcm.script"
$ws.Range("C57").Value2 = "scripts/"
$ws.Range("D57").Value2 = "Yes"
$ws.Range("E57").Value2 = "Utjan"

$ws.Range("A58").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B58").Value2 = "rax_icon_layers.script"
$ws.Range("C58").Value2 = "scripts/"
$ws.Range("D58").Value2 = "Yes"
$ws.Range("E58").Value2 = "Utjan"

$ws.Range("A59").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B59").Value2 = "z_item_icon_info.script"
$ws.Range("C59").Value2 = "scripts/"
$ws.Range("D59").Value2 = "Yes"
$ws.Range("E59").Value2 = "Utjan"

$ws.Range("A60").Value2 = "Utjans Item Ui Improvements"
$ws.Range("B60").Value2 = "Equipment_icons.dds"
$ws.Range("C60").Value2 = "textures/ui/"
$ws.Range("D60").Value2 = "Yes"
$ws.Range("E60").Value2 = "Utjan"

# ---------------------------------------------------------------------------
# 5. Update the view so the new rows are visible and A61 is selected, mimicking
#    the author's on-screen state when the workbook was saved.
# ---------------------------------------------------------------------------
$ws.Range("A61").Select() | Out-Null
